$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-07-04 Friday" "2025-07-05 Saturday"
Replace-Text "226÷4=" "781÷8="
Replace-Text "332÷4=" "541÷4="
Replace-Text "353÷6=" "774÷9="
Replace-Text "877÷8=" "987÷8="
Replace-Text "225÷5=" "542÷2="
Replace-Text "647÷2=" "291÷9="
Replace-Text "230÷2=" "488÷6="
Replace-Text "960÷2=" "474÷9="
Replace-Text "550÷7=" "934÷9="
Replace-Text "894÷2=" "604÷7="
Replace-Text "374÷4=" "158÷3="
Replace-Text "567÷7=" "850÷4="
Replace-Text "759÷7=" "211÷9="
Replace-Text "734÷6=" "576÷4="
Replace-Text "119÷8=" "945÷8="
Replace-Text "436÷8=" "842÷9="
Replace-Text "619÷3=" "462÷4="
Replace-Text "536÷6=" "267÷3="
Replace-Text "874÷2=" "353÷7="
Replace-Text "620÷2=" "486÷9="
Replace-Text "536÷2=" "713÷4="
Replace-Text "502÷6=" "465÷5="
Replace-Text "256÷2=" "882÷2="
Replace-Text "202÷9=" "663÷5="
Replace-Text "573÷6=" "614÷3="
